$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "City" column (A) now holds "City, State" header text, and a new
# column C is inserted that extracts just the city name from column A
# via a LEFT/SEARCH formula (split on the comma).
$ws.Range("A1").Value = "City, State"
$ws.Range("C1").Value = "City"

# C2 gets the formula typed directly; C3:C43 are filled down from it
# (R1C1 relative references so each row points at its own A cell).
$ws.Range("C2").Formula = "=LEFT(A2, SEARCH(`",`",A2,1)-1)"
$ws.Range("C3:C43").FormulaR1C1 = "=LEFT(RC[-2], SEARCH(`",`",RC[-2],1)-1)"

# Reflect the new selection / scroll position left by the edit.
[void]$ws.Range("C2:C43").Select()
